# Fixed index bug: the table was missing a row describing the column
# "types" (Number / String / String) right under the header row, which
# caused every data row below it to be off by one. Insert that row back
# in at row 2 - this naturally shifts every subsequent row (and its
# original values) down by one, matching the corrected data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the first data row (row 2), pushing the
# existing data rows down by one without altering their contents.
$ws.Rows.Item(2).EntireRow.Insert()

# Populate the newly inserted row with the missing "type" labels.
$ws.Range("A2").Value = "Number"
$ws.Range("B2").Value = "String"
$ws.Range("C2").Value = "String"
